# Adding averages and more checks
# - Refresh "PERIOD TO EXPIRE" (col H) and "LAST UPDATE" (col I) on the Training
#   Dashboard sheet: last update moved forward from 08-Sep-2025 to 16-Sep-2025,
#   which reduces every "days to expire" figure by 8.
# - Re-style the report headers: bold white text on the dark-blue header band
#   (and the title banner uses the same bold/white font, at the normal 11pt size
#   instead of the old 14pt).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Update the data on the "Training Dashboard" sheet (rows 3-20)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Training Dashboard")

# Force column I to stay plain text so the date string is not reinterpreted
# as a real date/serial number by Excel's automatic data-type detection.
$ws.Range("I3:I20").NumberFormat = "@"

for ($r = 3; $r -le 20; $r++) {
    $periodCell = $ws.Cells.Item($r, 8)
    $periodCell.Value2 = $periodCell.Value2 - 8

    $updateCell = $ws.Cells.Item($r, 9)
    $updateCell.Value = "16-Sep-2025"
}

# ---------------------------------------------------------------------------
# 2) Re-style the title banner and column-header rows on every sheet
# ---------------------------------------------------------------------------
foreach ($sheet in $wb.Worksheets) {
    $title = $sheet.Range("A1")
    $title.Font.Bold = $true
    $title.Font.Size = 11
    $title.Font.Color = 16777215

    $header = $sheet.Rows.Item(2)
    $header.Font.Bold = $true
    $header.Font.Color = 16777215
}
